$d = $word.ActiveDocument

$d.Content.Find.Execute("565÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷8=", 2) | Out-Null
$d.Content.Find.Execute("682÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "835÷6=", 2) | Out-Null
$d.Content.Find.Execute("311÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "498÷8=", 2) | Out-Null
$d.Content.Find.Execute("130÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "150÷3=", 2) | Out-Null
$d.Content.Find.Execute("521÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "589÷5=", 2) | Out-Null
$d.Content.Find.Execute("124÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "400÷9=", 2) | Out-Null
$d.Content.Find.Execute("161÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "123÷8=", 2) | Out-Null
$d.Content.Find.Execute("791÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "235÷5=", 2) | Out-Null
$d.Content.Find.Execute("385÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "310÷2=", 2) | Out-Null
$d.Content.Find.Execute("951÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "583÷7=", 2) | Out-Null
$d.Content.Find.Execute("652÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "516÷5=", 2) | Out-Null
$d.Content.Find.Execute("362÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "151÷6=", 2) | Out-Null
$d.Content.Find.Execute("541÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "482÷3=", 2) | Out-Null
$d.Content.Find.Execute("491÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "102÷5=", 2) | Out-Null
$d.Content.Find.Execute("572÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "420÷2=", 2) | Out-Null
$d.Content.Find.Execute("725÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "229÷6=", 2) | Out-Null
$d.Content.Find.Execute("293÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "415÷4=", 2) | Out-Null
$d.Content.Find.Execute("776÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "191÷9=", 2) | Out-Null
$d.Content.Find.Execute("854÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "207÷3=", 2) | Out-Null
$d.Content.Find.Execute("267÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "726÷5=", 2) | Out-Null
$d.Content.Find.Execute("100÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "839÷9=", 2) | Out-Null
$d.Content.Find.Execute("647÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "345÷7=", 2) | Out-Null
$d.Content.Find.Execute("767÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "562÷9=", 2) | Out-Null
$d.Content.Find.Execute("691÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "779÷2=", 2) | Out-Null
$d.Content.Find.Execute("240÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "631÷9=", 2) | Out-Null
